$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 48.77778
$ws.Range("I11").Value = 48.77778
$ws.Range("K11").Value = 48.77778
$ws.Range("M11").Value = 91.22221999999999
# Row 17
$ws.Range("H17").Value = 1580.4706
$ws.Range("J17").Value = 1580.4706
$ws.Range("L17").Value = 4741.4118
$ws.Range("N17").Value = -5077.4118
# Row 113
$ws.Range("H113").Value = 3908.1667
$ws.Range("I113").Value = 3250
$ws.Range("J113").Value = 4039.8
$ws.Range("K113").Value = 3250
$ws.Range("L113").Value = 4039.8
$ws.Range("M113").Value = 4
$ws.Range("N113").Value = -10547.8
# Row 135
$ws.Range("H135").Value = 446.66666
$ws.Range("I135").Value = 477.5
$ws.Range("J135").Value = 200
$ws.Range("K135").Value = 4297.5
$ws.Range("L135").Value = 1800
$ws.Range("M135").Value = -1762.5
$ws.Range("N135").Value = -6870
# Row 137
$ws.Range("H137").Value = 1280.75
$ws.Range("I137").Value = 1041.1666
$ws.Range("K137").Value = 3123.4998
$ws.Range("M137").Value = -573.4998000000001
# Row 141
$ws.Range("H141").Value = 3544.4
$ws.Range("I141").Value = 3492.4285
$ws.Range("J141").Value = 3665.6667
$ws.Range("K141").Value = 10477.2855
$ws.Range("L141").Value = 10997.0001
$ws.Range("M141").Value = -5297.2855
$ws.Range("N141").Value = -21357.0001

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 352.83334
$ws.Range("I2").Value = 352.83334
$ws.Range("K2").Value = 352.83334
$ws.Range("M2").Value = -239.83334
# Row 74
$ws.Range("H74").Value = 694.75
$ws.Range("I74").Value = 643.5
$ws.Range("J74").Value = 848.5
$ws.Range("K74").Value = 643.5
$ws.Range("L74").Value = 848.5
$ws.Range("M74").Value = 230.5
$ws.Range("N74").Value = -2596.5
# Row 77
$ws.Range("H77").Value = 694.75
$ws.Range("I77").Value = 643.5
$ws.Range("J77").Value = 848.5
$ws.Range("K77").Value = 3217.5
$ws.Range("L77").Value = 4242.5
$ws.Range("M77").Value = 1150.5
$ws.Range("N77").Value = -12978.5
# Row 116
$ws.Range("H116").Value = 352.83334
$ws.Range("I116").Value = 352.83334
$ws.Range("K116").Value = 352.83334
$ws.Range("M116").Value = 1941.16666
# Row 132
$ws.Range("H132").Value = 2902.6667
$ws.Range("I132").Value = 2902.6667
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8708.000100000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6178.000100000001
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 352.83334
$ws.Range("I3").Value = 352.83334
$ws.Range("K3").Value = 352.83334
$ws.Range("M3").Value = -238.83334
# Row 81
$ws.Range("H81").Value = 35629.668
$ws.Range("J81").Value = 35629.668
$ws.Range("L81").Value = 35629.668
$ws.Range("N81").Value = -37751.668
# Row 84
$ws.Range("H84").Value = 35629.668
$ws.Range("J84").Value = 35629.668
$ws.Range("L84").Value = 106889.004
$ws.Range("N84").Value = -117497.004
# Row 94
$ws.Range("H94").Value = 2500
$ws.Range("I94").Value = 2500
$ws.Range("K94").Value = 2500
$ws.Range("M94").Value = -2049
# Row 138
$ws.Range("H138").Value = 99997.5
$ws.Range("J138").Value = 99997.5
$ws.Range("L138").Value = 99997.5
$ws.Range("N138").Value = -110277.5

$ws = $wb.Worksheets.Item("CRP")
# Row 23
$ws.Range("H23").Value = 250000
$ws.Range("I23").Value = 250000
$ws.Range("K23").Value = 250000
$ws.Range("M23").Value = -249760
# Row 27
$ws.Range("H27").Value = 250000
$ws.Range("I27").Value = 250000
$ws.Range("K27").Value = 250000
$ws.Range("M27").Value = -249808
# Row 31
$ws.Range("H31").Value = 979.2
$ws.Range("I31").Value = 974.75
$ws.Range("K31").Value = 974.75
$ws.Range("M31").Value = -679.75
# Row 34
$ws.Range("H34").Value = 979.2
$ws.Range("I34").Value = 974.75
$ws.Range("K34").Value = 974.75
$ws.Range("M34").Value = -772.75
# Row 41
$ws.Range("H41").Value = 1000
$ws.Range("I41").Value = 1000
$ws.Range("K41").Value = 1000
$ws.Range("M41").Value = -572
# Row 47
$ws.Range("H47").Value = 20000
$ws.Range("I47").Value = 20000
$ws.Range("K47").Value = 20000
$ws.Range("M47").Value = -19434
# Row 134
$ws.Range("H134").Value = 3001.125
$ws.Range("I134").Value = 3001.125
$ws.Range("K134").Value = 9003.375
$ws.Range("M134").Value = -6468.375

$ws = $wb.Worksheets.Item("CUL")
# Row 51
$ws.Range("H51").Value = 500
$ws.Range("I51").Value = 500
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 1500
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -1040
$ws.Range("N51").ClearContents()
# Row 80
$ws.Range("H80").Value = 8999
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 8999
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 26997
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -28869
# Row 83
$ws.Range("H83").Value = 8999
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 8999
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 80991
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -90351
# Row 92
$ws.Range("H92").Value = 439
$ws.Range("I92").Value = 198.66667
$ws.Range("J92").Value = 645
$ws.Range("K92").Value = 596.00001
$ws.Range("L92").Value = 1935
$ws.Range("M92").Value = 651.99999
$ws.Range("N92").Value = -4431

$ws = $wb.Worksheets.Item("GSM")
# Row 140
$ws.Range("H140").Value = 142948.5
$ws.Range("J140").Value = 142948.5
$ws.Range("L140").Value = 142948.5
$ws.Range("N140").Value = -153308.5

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
# Row 71
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
# Row 93
$ws.Range("H93").Value = 3081.5
$ws.Range("I93").Value = 3081.5
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 3081.5
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -1833.5
$ws.Range("N93").ClearContents()
# Row 122
$ws.Range("H122").Value = 6603.125
$ws.Range("I122").Value = 5041.4165
$ws.Range("J122").Value = 8164.8335
$ws.Range("K122").Value = 15124.2495
$ws.Range("L122").Value = 24494.5005
$ws.Range("M122").Value = -12674.2495
$ws.Range("N122").Value = -29394.5005

$ws = $wb.Worksheets.Item("WVR")
# Row 39
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
# Row 81
$ws.Range("H81").Value = 10000002
$ws.Range("J81").Value = 10000002
$ws.Range("L81").Value = 20000004
$ws.Range("N81").Value = -20002126
# Row 84
$ws.Range("H84").Value = 10000002
$ws.Range("J84").Value = 10000002
$ws.Range("L84").Value = 100000020
$ws.Range("N84").Value = -100010628
# Row 101
$ws.Range("H101").Value = 16151
$ws.Range("J101").Value = 16151
$ws.Range("L101").Value = 16151
$ws.Range("N101").Value = -22641
# Row 132
$ws.Range("H132").Value = 3009.8
$ws.Range("I132").Value = 3009.8
$ws.Range("K132").Value = 9029.400000000001
$ws.Range("M132").Value = -6499.400000000001
